$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell for the "password" column (N)
$ws.Range("N1").Value = "password"

# --- Row 23: Yash ---
$ws.Range("B23").Value = "Yash"
$ws.Range("C23").Value = "BCA"
$ws.Range("D23").Value = "yash"
$ws.Range("E23").Value = "R"
$ws.Range("F23").Value = "prajapati"
$ws.Range("G23").Value = "bholad"
$ws.Range("H23").Value = 9638845196
$ws.Range("J23").Value = "yash@gmail.com"
$ws.Range("K23").Value = "M"
$ws.Range("L23").Value = "2006-04-11"
$ws.Range("M23").Value = "photos/2029-008.jpg"
$ws.Range("A23").Value = "2019-008"

# --- Row 24: Mukesh ---
$ws.Range("A24").Value = "2022-037"
$ws.Range("B24").Value = "Mukesh"
$ws.Range("C24").Value = 11
$ws.Range("D24").Value = "Mukesh"
$ws.Range("E24").Value = "Mahadevbhai"
$ws.Range("F24").Value = "kolipatel"
$ws.Range("G24").Value = "Kamijala"
$ws.Range("H24").Value = 8347898768
$ws.Range("J24").Value = "mukesh@gmail.com"
$ws.Range("K24").Value = "Male"
$ws.Range("L24").Value = "2009-04-24"
$ws.Range("M24").Value = "photos/2022-037.jpg"

# Hyperlinks (mailto) for the new email cells, matching J22's existing pattern
$ws.Hyperlinks.Add($ws.Range("J23"), "mailto:yash@gmail.com")
$ws.Hyperlinks.Add($ws.Range("J24"), "mailto:mukesh@gmail.com")

# Re-apply the Hyperlink cell style so J23/J24 match J22's look exactly
$ws.Range("J23").Style = "Hyperlink"
$ws.Range("J24").Style = "Hyperlink"

# Column widths settle slightly as Excel recalculates "best fit" for the
# wider photo-path text and the new "password" column
$ws.Columns.Item(13).ColumnWidth = 18.42
$ws.Columns.Item(14).ColumnWidth = 8.6

# View state: zoom to 100% and move the selection near the new rows
$win = $excel.Windows.Item(1)
$win.Zoom = 100
$ws.Range("J30").Select()
